$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.969.46'
$ws.Range('E2').Value = '  -2.70%  '
$ws.Range('D3').Value = '3.366.96'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '567.85'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.35'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.97'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.417'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.62%  '
$ws.Range('D12').Value = '3.950.63'
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.00'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').Value = '3.378.20'
$ws.Range('E15').Value = '  -2.25%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000169'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '61.018.23'
$ws.Range('E17').Value = '  -2.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.34'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '14.46'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '375.35'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.07%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.563'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '75.47'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '3.506.24'
$ws.Range('E25').Value = '  -2.36%  '
$ws.Range('E26').Value = '  -6.65%  '
$ws.Range('E27').Value = '  -3.58%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.44'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.72'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.49%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '22.86'
$ws.Range('D33').Style = "Normal"
$ws.Range('E34').Value = '  -3.96%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.36'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '169.85'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.82'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '29.15'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -8.36%  '
$ws.Range('D40').Value = '3.402.70'
$ws.Range('E40').Value = '  -2.14%  '
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.761'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.42%  '
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.14'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.77%  '
$ws.Range('E45').Value = '  -4.84%  '
$ws.Range('D46').Value = '2.494.87'
$ws.Range('E46').Value = '  -2.47%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '22.81'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '6.68'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.92%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0261'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.819'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.59%  '
